# Update cryptos list: refresh Price/Volume(1h) columns, and restore
# the two coin-rank swaps (Stacks<->EnergySwap, InjectiveProtocol<->Bittensor).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.659.92"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "3.188.14"
$ws.Range("E3").Value = "  -4.56%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'571.69"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("D6").Value = "'171.85"
$ws.Range("E6").Value = "  -2.96%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'0.601"
$ws.Range("E8").Value = "  -3.14%  "
$ws.Range("D9").Value = "3.192.44"
$ws.Range("E9").Value = "  -4.34%  "
$ws.Range("D10").Value = "'0.125"
$ws.Range("E10").Value = "  -3.26%  "
$ws.Range("D11").Value = "'6.56"
$ws.Range("E11").Value = "  -4.66%  "
$ws.Range("D12").Value = "'0.392"
$ws.Range("E12").Value = "  -4.50%  "
$ws.Range("D13").Value = "3.737.55"
$ws.Range("E13").Value = "  -4.77%  "
$ws.Range("E14").Value = "  +1.49%  "
$ws.Range("D15").Value = "'27.35"
$ws.Range("E15").Value = "  -4.73%  "
$ws.Range("D16").Value = "65.674.54"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").Value = "'0.0000164"
$ws.Range("E17").Value = "  -3.37%  "
$ws.Range("D18").Value = "3.176.98"
$ws.Range("E18").Value = "  -4.40%  "
$ws.Range("D19").Value = "'5.72"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "'12.87"
$ws.Range("E20").Value = "  -4.26%  "
$ws.Range("D21").Value = "'359.63"
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("D22").Value = "'7.28"
$ws.Range("E22").Value = "  -2.18%  "
$ws.Range("E23").Value = "  +0.76%  "
$ws.Range("D24").Value = "'69.21"
$ws.Range("E24").Value = "  -3.55%  "
$ws.Range("D25").Value = "'0.495"
$ws.Range("E25").Value = "  -4.97%  "
$ws.Range("D26").Value = "3.288.96"
$ws.Range("E26").Value = "  -6.22%  "
$ws.Range("D27").Value = "'0.0000116"
$ws.Range("E27").Value = "  -5.29%  "
$ws.Range("D28").Value = "'9.81"
$ws.Range("E28").Value = "  +2.14%  "
$ws.Range("E29").Value = "  -1.24%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  -0.33%  "
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("D32").Value = "'1.93"
$ws.Range("E32").Value = "  -2.07%  "
$ws.Range("D33").Value = "'5.37"
$ws.Range("E33").Value = "  -5.29%  "
$ws.Range("D34").Value = "'21.98"
$ws.Range("E34").Value = "  -3.64%  "
$ws.Range("E35").Value = "  -1.77%  "
$ws.Range("D36").Value = "'6.61"
$ws.Range("E36").Value = "  -3.91%  "
$ws.Range("D37").Value = "'159.95"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("E38").Value = "  -3.40%  "
$ws.Range("D39").Value = "'0.835"
$ws.Range("E39").Value = "  -1.70%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").Value = "'26.50"
$ws.Range("E40").Value = "  -2.52%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.79"
$ws.Range("E41").Value = "  +1.91%  "
$ws.Range("D42").Value = "'2.50"
$ws.Range("E42").Value = "  -1.86%  "
$ws.Range("D43").Value = "2.647.90"
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("E44").Value = "  -2.03%  "
$ws.Range("E45").Value = "  -2.04%  "
$ws.Range("D46").Value = "'39.61"
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("D47").Value = "'0.0660"
$ws.Range("E47").Value = "  -0.83%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").Value = "'329.43"
$ws.Range("E48").Value = "  -1.93%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "'24.07"
$ws.Range("E49").Value = "  -1.41%  "
$ws.Range("E50").Value = "  -1.62%  "
$ws.Range("E51").Value = "  -1.09%  "
